# Apply the PIB_SCT.xlsx update:
#  - Insert a new row for 2024 Q4 data at the top of the data table (row 5),
#    shifting the existing rows (and footnotes) down by one.
#  - Refresh the revised figures for 2024 Q3/Q2/Q1 (now rows 6-8).
#  - Update the "Actualización" note text/date (now row 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_1.3")

# Insert a new row above the current row 5 (2024 Q3), pushing everything
# below (including the footnote rows) down by one.
$ws.Rows("5:5").Insert()

# New row 5: 2024 Q4
$ws.Range("B5").Value = 2024
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 1890306.1189999999
$ws.Range("E5").Value = 86251.111999999994
$ws.Range("F5").Value = 44629.324999999997
$ws.Range("G5").Value = 1436991.175
$ws.Range("H5").Value = 37502.938999999998
$ws.Range("I5").Value = 478008.83100000001
$ws.Range("J5").Value = 66976.903000000006
$ws.Range("K5").Value = 360880.18599999999

# Copy the formatting from the row two below (2024 Q2, now row 7, which
# carries the "even" row style) so the new row matches the alternating
# banding used throughout the table.
$ws.Range("B7:K7").Copy()
$ws.Range("B5:K5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Revised figures for 2024 Q3 (row 6)
$ws.Range("D6").Value = 1899914.575
$ws.Range("E6").Value = 84544.587
$ws.Range("G6").Value = 1459757.4959999998
$ws.Range("I6").Value = 453199.19300000003
$ws.Range("J6").Value = 57875.857000000004
$ws.Range("K6").Value = 341997.723

# Revised figures for 2024 Q2 (row 7)
$ws.Range("D7").Value = 1939608.1529999999
$ws.Range("G7").Value = 1502603.7519999999

# Revised figures for 2024 Q1 (row 8)
$ws.Range("D8").Value = 1830286.128
$ws.Range("G8").Value = 1420126.96

# Update the "Actualización" note (now on row 33 after the insert).
$ws.Range("B33").Value = "Actualización: Enero 2025."
